$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Rows 1-4: swap the leading summary values for placeholder values ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "366"

# --- Row 6: 0.00004 -> 0.00006 ---
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00006"

# --- Remove row 9 (0.00002) entirely ---
$t.Rows.Item(9).Delete()

# After the delete, the old row10/row11/row12 shifted up to row9/row10/row11.
$t.Rows.Item(9).Cells.Item(1).Range.Text = "0.00004"
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00004"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00004"

# --- Insert a brand-new row (0.01138) right before the old row12 (100.0) ---
$beforeRow = $t.Rows.Item(12)
$t.Rows.Add($beforeRow) | Out-Null
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.01138"

# --- Collapse the final three multi-value (tab-separated) rows down to a single value each ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.99"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.01"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "179"
